$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.896.22'
$ws.Range('E2').Value = '  +0.03%  '
$ws.Range('D3').Value = '1.731.34'
$ws.Range('E3').Value = '  -0.53%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '245.36'
$ws.Range('E5').Value = '  +2.97%  '
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('D7').Value = '0.5031'
$ws.Range('E7').Value = '  -2.70%  '
$ws.Range('D8').Value = '0.2717'
$ws.Range('E8').Value = '  -1.08%  '
$ws.Range('D9').Value = '0.06169'
$ws.Range('E9').Value = '  +0.20%  '
$ws.Range('D10').Value = '1.732.39'
$ws.Range('E10').Value = '  -0.47%  '
$ws.Range('E11').Value = '  +0.96%  '
$ws.Range('D12').Value = '0.6531'
$ws.Range('E12').Value = '  +1.32%  '
$ws.Range('D13').Value = '15.19'
$ws.Range('E13').Value = '  +1.49%  '
$ws.Range('D14').Value = '4.757'
$ws.Range('E14').Value = '  +3.44%  '
$ws.Range('D15').Value = '76.97'
$ws.Range('E15').Value = '  -0.60%  '
$ws.Range('D16').Value = '1.000'
$ws.Range('E16').Value = '  +0.00%  '
$ws.Range('D17').Value = '1.000'
$ws.Range('E17').Value = '  -0.01%  '
$ws.Range('D18').Value = '25.902.77'
$ws.Range('E18').Value = '  +0.01%  '
$ws.Range('D19').Value = '11.90'
$ws.Range('E19').Value = '  +1.58%  '
$ws.Range('D20').Value = '0.000006825'
$ws.Range('E20').Value = '  +0.82%  '
$ws.Range('D21').Value = '4.589'
$ws.Range('E21').Value = '  +7.36%  '
$ws.Range('D22').Value = '1.958.25'
$ws.Range('E22').Value = '  -0.33%  '
$ws.Range('D23').Value = '8.805'
$ws.Range('E23').Value = '  +1.70%  '
$ws.Range('D24').Value = '5.479'
$ws.Range('E24').Value = '  +4.43%  '
$ws.Range('D25').Value = '134.49'
$ws.Range('E25').Value = '  -3.07%  '
$ws.Range('E26').Value = '  +1.14%  '
$ws.Range('D27').Value = '1.421'
$ws.Range('E27').Value = '  -5.73%  '
$ws.Range('D28').Value = '1.788'
$ws.Range('E28').Value = '  +1.67%  '
$ws.Range('D29').Value = '105.44'
$ws.Range('E29').Value = '  -0.42%  '
$ws.Range('D30').Value = '3.963'
$ws.Range('E30').Value = '  +0.47%  '
$ws.Range('D31').Value = '0.08130'
$ws.Range('E31').Value = '  -1.90%  '
$ws.Range('D32').Value = '3.700'
$ws.Range('E32').Value = '  +0.67%  '
$ws.Range('D33').Value = '0.04723'
$ws.Range('E33').Value = '  +2.78%  '
$ws.Range('D34').Value = '2.651'
$ws.Range('E34').Value = '  +0.28%  '
$ws.Range('D35').Value = '0.9945'
$ws.Range('E35').Value = '  +0.75%  '
$ws.Range('D36').Value = '0.6133'
$ws.Range('E36').Value = '  -0.78%  '
$ws.Range('D37').Value = '2.746'
$ws.Range('E37').Value = '  +2.50%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '0.01610'
$ws.Range('E38').Value = '  +0.01%  '
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').Value = '0.8800'
$ws.Range('E39').Value = '  +18.81%  '
$ws.Range('D40').Value = '1.960'
$ws.Range('E40').Value = '  +1.89%  '
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('D42').Value = '101.78'
$ws.Range('E42').Value = '  +4.08%  '
$ws.Range('D43').Value = '0.3907'
$ws.Range('E43').Value = '  +1.98%  '
$ws.Range('E44').Value = '  +0.53%  '
$ws.Range('D45').Value = '0.1184'
$ws.Range('E45').Value = '  +5.03%  '
$ws.Range('D46').Value = '6.367'
$ws.Range('E46').Value = '  +2.73%  '
$ws.Range('D47').Value = '55.72'
$ws.Range('E47').Value = '  +1.36%  '
$ws.Range('E48').Value = '  +0.35%  '
$ws.Range('D49').Value = '30.79'
$ws.Range('E49').Value = '  +1.08%  '
$ws.Range('B50').Value = 'Decentraland'
$ws.Range('C50').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D50').Value = '0.3488'
$ws.Range('E50').Value = '  +2.56%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = '7.634'
$ws.Range('E51').Value = '  +0.18%  '
